# Stage 6. Added an opportunity to change date of group creation
# & check for uniqueness of roles' names.

$wb = $excel.ActiveWorkbook

# --- Sheet "Splin": add a "Date of group creation" column (H) ---
$splin = $wb.Worksheets.Item("Splin")

# Header cell H1 - copy the header formatting (fill/style) from G1, then set the text.
$splin.Range("G1").Copy()
$splin.Range("H1").PasteSpecial(-4122)
$splin.Range("H1").Value = "Date of group creation"

# Widen column H so the new header text fits.
$splin.Columns.Item(8).ColumnWidth = 19

# H4 - the date the group was created, formatted like the other date cells (copy from C2).
$splin.Range("C2").Copy()
$splin.Range("H4").PasteSpecial(-4122)
$splin.Range("H4").Value = [DateTime]"1994-05-27"

$splin.Range("D12").Select()

# --- Sheet "Iriao": remove the duplicate-named row (row 4) ---
$iriao = $wb.Worksheets.Item("Iriao")

$iriao.Range("A4:G4").ClearContents()

$iriao.Range("J6").Select()
